$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column BB: next quarterly date in the header row, plus one more
# column of forecast data appended to every data row (rows 3-18 simply
# carry forward the plateaued BA value, rows 19-21 get new forecast
# values).

# Row 1 header date (copy style/format from BA1, which already uses the
# worksheet's date style).
$ws.Range("BB1").Value = 45986
$ws.Range("BA1").Copy()
$ws.Range("BB1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Rows 3-18: repeat the last (BA) value into the new BB column.
$ws.Range("BB3").Value = 0.2361821772587591
$ws.Range("BB4").Value = 0.1732386938003039
$ws.Range("BB5").Value = 2.130327852244007
$ws.Range("BB6").Value = 1.019364148315116
$ws.Range("BB7").Value = 0.8966256141480056
$ws.Range("BB8").Value = 0.4602645177979703
$ws.Range("BB9").Value = 2.152035874345892
$ws.Range("BB10").Value = 1.918040486507278
$ws.Range("BB11").Value = 1.673405723817356
$ws.Range("BB12").Value = 0.7965036850253515
$ws.Range("BB13").Value = 1.420403924237745
$ws.Range("BB14").Value = -4.401478753282307
$ws.Range("BB15").Value = -1.740831309918522
$ws.Range("BB16").Value = 5.27347369917277
$ws.Range("BB17").Value = -0.7058256052666523
$ws.Range("BB18").Value = 0.4368922194300628

# Rows 19-21: new forecast values (not a straight carry-forward).
$ws.Range("BB19").Value = 0.8976398032236155
$ws.Range("BB20").Value = 0.4275768375374467
$ws.Range("BB21").Value = 0.6589725835419058
